$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.816.57'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '3.390.78'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.20'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.71'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.55%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '3.389.83'
$ws.Range('E8').Value = '  -1.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.474'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.53'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.125'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.394'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').Value = '3.969.52'
$ws.Range('E13').Value = '  -1.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.29'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.124'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('E16').Value = '  -2.38%  '
$ws.Range('D17').Value = '3.395.50'
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('D18').Value = '60.953.15'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.31'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.22'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.12'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.56'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.565'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.08'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('E26').Value = '  -3.18%  '
$ws.Range('D27').Value = '3.533.93'
$ws.Range('E27').Value = '  -1.43%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.35'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -6.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.19'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.16'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.44'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -7.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.82'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('E36').Value = '  -0.52%  '
$ws.Range('D37').Value = '3.417.37'
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.08'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.24%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '167.22'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.54'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0783'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '27.16'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.10%  '
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.50'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.87'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.70'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.28%  '
$ws.Range('D48').Value = '2.560.19'
$ws.Range('E48').Value = '  -1.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.13'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.90'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.89'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.70%  '
